# Staff List: insert Username/Password columns, lowercase role values, re-layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing 5x4 data block (Staff ID, Name, Role, Gender, Age) before
# we shift everything over by inserting the two new leading columns.
$oldHeaders = @()
for ($c = 1; $c -le 5; $c++) { $oldHeaders += $ws.Cells.Item(1, $c).Value2 }

$oldData = @()
for ($r = 2; $r -le 5; $r++) {
    $row = @()
    for ($c = 1; $c -le 5; $c++) { $row += $ws.Cells.Item($r, $c).Value2 }
    $oldData += ,$row
}

# Insert two new columns at the front for Username / Password.
$ws.Range("A1:B1").EntireColumn.Insert()

# Usernames / passwords / lowercase roles for each data row.
$usernames = @("doctor1", "doctor2", "pharmacist1", "admin1")
$passwords = @("password", "password", "password", "password")
$roles = @("doctor", "doctor", "pharmacist", "administrator")

# Write column-by-column (matches the shared-string insertion order of the
# source edit: Username col, then Password col, then the shifted columns).
$ws.Range("A1").Value2 = "Username"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 1).Value2 = $usernames[$i] }

$ws.Range("B1").Value2 = "Password"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 2).Value2 = $passwords[$i] }

$ws.Range("C1").Value2 = "Staff ID"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 3).Value2 = $oldData[$i][0] }

$ws.Range("D1").Value2 = "Name"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 4).Value2 = $oldData[$i][1] }

$ws.Range("E1").Value2 = "Role"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 5).Value2 = $roles[$i] }

$ws.Range("F1").Value2 = "Gender"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 6).Value2 = $oldData[$i][3] }

$ws.Range("G1").Value2 = "Age"
for ($i = 0; $i -lt 4; $i++) { $ws.Cells.Item($i + 2, 7).Value2 = $oldData[$i][4] }

# The two new header cells need the same bold/centered/bordered header style
# as the rest of row 1.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Column widths (values picked so that this host's char-width -> pixel ->
# char-width rounding lands on the same stored width as the target file).
$ws.Range("A1").EntireColumn.ColumnWidth = 10.25
$ws.Range("B1").EntireColumn.ColumnWidth = 10.25
$ws.Range("D1").EntireColumn.ColumnWidth = 24
$ws.Range("E1").EntireColumn.ColumnWidth = 29.65

# Selection / used range.
$ws.Range("C3").Select() | Out-Null
